$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely; this shifts all rows below it up by one,
# effectively removing the stale "Lunes 03/06/2024" entry and shrinking
# the used range from A1:B7 to A1:B6.
$ws.Rows("2:2").Delete()
